# Add two new Mac-Addresses (10 new rows: regcntr 10001 x machines 10030/10031, each with 5 new devices)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newRows = @(
    @(147, 10001, 10030, 3000166),
    @(148, 10001, 10030, 3000167),
    @(149, 10001, 10030, 3000168),
    @(150, 10001, 10030, 3000169),
    @(151, 10001, 10030, 3000170),
    @(152, 10001, 10031, 3000171),
    @(153, 10001, 10031, 3000172),
    @(154, 10001, 10031, 3000173),
    @(155, 10001, 10031, 3000174),
    @(156, 10001, 10031, 3000175)
)

foreach ($row in $newRows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = "eng"
    $ws.Cells.Item($r, 5).Value = $true
    $ws.Cells.Item($r, 6).Value = "superadmin"
    $ws.Cells.Item($r, 7).Value = "now()"
}

$win = $excel.ActiveWindow
$win.ScrollRow = 144
$win.ScrollColumn = 1
[void]$ws.Range("H149").Select()
